$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.709.00'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '2.474.95'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''319.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D6").Value = '''93.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.26%  '
$ws.Range("D7").Value = '''0.552'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.57%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '''0.517'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("B10").Value = 'Avalanche'
$ws.Range("C10").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D10").Value = '''33.29'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.10%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '''0.0864'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.41%  '
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").Value = '2.859.12'
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").Value = '''15.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.35%  '
$ws.Range("D16").Value = '2.479.53'
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("E17").Value = '  +2.37%  '
$ws.Range("D18").Value = '41.707.26'
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D19").Value = '''6.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '0.0₃0951'
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("D21").Value = '''71.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").Value = '''11.29'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.02%  '
$ws.Range("D23").Value = '''239.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.36%  '
$ws.Range("E24").Value = '  +1.39%  '
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Value = '''24.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").Value = '''9.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.25%  '
$ws.Range("D30").Value = '''36.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.22%  '
$ws.Range("D31").Value = '''158.62'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("D32").Value = '''5.51'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("E35").Value = '  +1.17%  '
$ws.Range("D36").Value = '''17.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.77%  '
$ws.Range("E37").Value = '  +5.85%  '
$ws.Range("E38").Value = '  +2.00%  '
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("D41").Value = '''4.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.88%  '
$ws.Range("E42").Value = '  +11.08%  '
$ws.Range("D43").Value = '1.993.67'
$ws.Range("E43").Value = '  +2.47%  '
$ws.Range("D44").Value = '''0.0285'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("D45").Value = '''18.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.90%  '
$ws.Range("E46").Value = '  +2.42%  '
$ws.Range("D47").Value = '''9.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.69%  '
$ws.Range("D48").Value = '2.716.55'
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("D49").Value = '''97.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.68%  '
$ws.Range("D50").Value = '''73.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.36%  '
$ws.Range("D51").Value = '''67.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.14%  '
